# tracker.xlsx update
# - Add "Human" (affected species) to row 15 (J15)
# - Turn P15's ProMED link text into a real hyperlink (adds Hyperlink cell style)
# - Append a new row (43) for a Dengue / Samoa ProMED alert
# - Update the sheet's scroll/selection state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: fill in the missing "Affected species" (Human) ---
$ws.Range("J15").Value = "Human"

# --- Row 15: make the ProMED link in P15 a real hyperlink ---
[void]$ws.Hyperlinks.Add($ws.Range("P15"), "https://promedmail.org/promed-post/?id=8716004")

# --- New row 43: Dengue alert in Samoa ---
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "Dengue"
$ws.Range("E43").Value = "Samoa"
$ws.Range("F43").Value = -13.769389500000001
$ws.Range("G43").Value = -172.12004999999999
$ws.Range("H43").Value = "Virus"
$ws.Range("I43").Value = "Dengue Virus (DENV)"
$ws.Range("J43").Value = "Human"
$ws.Range("L43").Value = 45405
$ws.Range("L43").NumberFormat = "d-mmm-yy"
$ws.Range("P43").Value = "https://promedmail.org/promed-post/?id=8716106"

# --- Sheet view: keep the header row frozen, scroll down and reselect ---
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
[void]($win.FreezePanes = $true)
$ws.Range("H19").Select()
